$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price reading (Fecha 2021-12-07 / serial 44537) was added as
# the newest record. It goes in as the new row 2, pushing the three existing
# data rows (old rows 2-4) down to rows 3-5, unchanged.
$ws.Rows("2:2").Insert()

# Insert() copies formatting down from the header row (row 1, bold + border),
# which is not what the other data rows look like. Re-apply the plain data
# row formatting (incl. the date number format on column D) by copying it
# from the row immediately below (the former row 2, now row 3).
$ws.Range("A3:R3").Copy()
$ws.Range("A2:R2").PasteSpecial(-4122)

# Populate the new row with the new weekly reading.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C2").Value = "Arica y Parinacota"
$ws.Range("D2").Value = 44537
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112030
$ws.Range("G2").Value = "Poroto granado"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 1300
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = 1350
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 1350
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
